# "Removed table in course student"
#
# The EmailKey/Fullname roster on Sheet1 gains a new student
# (Saicharan Gurudu) and Sravya Kancharla's row -- which used to sit at
# the bottom of the sheet with a "s@nwmissouri.edu" mailto hyperlink --
# is corrected to use her real email and moved up under Girish Guntuku.
# The old hyperlinked placeholder row is dropped entirely, along with the
# now-unused built-in "Hyperlink" cell style, and the view settings are
# reset (no more frozen/scrolled top-left cell, new selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the two new rows right after "Girish Guntuku" (row 7).
$ws.Rows("8:9").Insert()

# New student record.
$ws.Range("A8").Value = "s531499@nwmissouri.edu"
$ws.Range("B8").Value = "Saicharan Gurudu"

# Sravya Kancharla, now with her correct email address, in her new slot.
$ws.Range("A9").Value = "s531500@nwmissouri.edu"
$ws.Range("B9").Value = "Sravya Kancharla"

# Drop the mailto: hyperlink that lived on the old trailing row, then
# remove that now-duplicated row (it shifted down to row 23 after the
# insert above).
$null = $ws.Hyperlinks.Delete()
$ws.Rows("23").Delete()

# The workbook no longer needs the built-in "Hyperlink" cell style since
# nothing uses it any more.
$wb.Styles("Hyperlink").Delete()

# Reset the view: clear the scrolled-down top-left cell and move the
# selection.
$null = $ws.Range("J13").Select()
